# "Répartition des taches" - reassign task owners
#   Row 5  (A5 = "10: skills")               -> was Roméo, now Fabio
#   Row 11 (A11 = "tache 14: équipement")     -> was Fabio, now Roméo
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

$ws.Range("B5").Value = "Fabio"
$ws.Range("B11").Value = "Roméo"

# Leave the selection on the cell that was last edited
$ws.Range("B11").Select()
